# Apply updated "想去人数" (want-to-go count, column F) values to match
# the latest scraped data (gh-pages output regenerated at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1572
$ws.Cells.Item(3, 6).Value = 3318
$ws.Cells.Item(4, 6).Value = 27
$ws.Cells.Item(5, 6).Value = 749
$ws.Cells.Item(6, 6).Value = 2345
$ws.Cells.Item(7, 6).Value = 506
$ws.Cells.Item(8, 6).Value = 423
$ws.Cells.Item(9, 6).Value = 253
$ws.Cells.Item(10, 6).Value = 148
$ws.Cells.Item(11, 6).Value = 369
$ws.Cells.Item(12, 6).Value = 1110
$ws.Cells.Item(13, 6).Value = 463
$ws.Cells.Item(16, 6).Value = 269
$ws.Cells.Item(17, 6).Value = 4858
$ws.Cells.Item(18, 6).Value = 27
$ws.Cells.Item(19, 6).Value = 1376
$ws.Cells.Item(20, 6).Value = 3568
$ws.Cells.Item(21, 6).Value = 142
$ws.Cells.Item(22, 6).Value = 201
$ws.Cells.Item(23, 6).Value = 3865
$ws.Cells.Item(24, 6).Value = 5233
$ws.Cells.Item(25, 6).Value = 126
$ws.Cells.Item(26, 6).Value = 986
$ws.Cells.Item(27, 6).Value = 574
$ws.Cells.Item(28, 6).Value = 3365
$ws.Cells.Item(29, 6).Value = 386
$ws.Cells.Item(30, 6).Value = 53
$ws.Cells.Item(31, 6).Value = 148
$ws.Cells.Item(32, 6).Value = 98
$ws.Cells.Item(33, 6).Value = 899
$ws.Cells.Item(34, 6).Value = 1217
$ws.Cells.Item(35, 6).Value = 31
$ws.Cells.Item(36, 6).Value = 47
$ws.Cells.Item(37, 6).Value = 1440
$ws.Cells.Item(38, 6).Value = 145
$ws.Cells.Item(39, 6).Value = 1424
$ws.Cells.Item(40, 6).Value = 31
$ws.Cells.Item(41, 6).Value = 916
$ws.Cells.Item(42, 6).Value = 892
$ws.Cells.Item(43, 6).Value = 525
$ws.Cells.Item(45, 6).Value = 1284
$ws.Cells.Item(46, 6).Value = 87
$ws.Cells.Item(47, 6).Value = 184
$ws.Cells.Item(49, 6).Value = 3757

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 1026
$ws.Cells.Item(14, 6).Value = 6
$ws.Cells.Item(15, 6).Value = 8

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 2445

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 2445
$ws.Cells.Item(3, 6).Value = 1572
$ws.Cells.Item(4, 6).Value = 3318
$ws.Cells.Item(5, 6).Value = 27
$ws.Cells.Item(6, 6).Value = 749
$ws.Cells.Item(8, 6).Value = 2345
$ws.Cells.Item(9, 6).Value = 506
$ws.Cells.Item(10, 6).Value = 423
$ws.Cells.Item(11, 6).Value = 253
$ws.Cells.Item(12, 6).Value = 1026
$ws.Cells.Item(13, 6).Value = 148
$ws.Cells.Item(14, 6).Value = 369
$ws.Cells.Item(15, 6).Value = 1110
$ws.Cells.Item(16, 6).Value = 463
$ws.Cells.Item(19, 6).Value = 269
$ws.Cells.Item(20, 6).Value = 4858
$ws.Cells.Item(22, 6).Value = 1376
$ws.Cells.Item(23, 6).Value = 3865
$ws.Cells.Item(24, 6).Value = 5233
$ws.Cells.Item(25, 6).Value = 126
$ws.Cells.Item(26, 6).Value = 986
$ws.Cells.Item(27, 6).Value = 574
$ws.Cells.Item(28, 6).Value = 3365
$ws.Cells.Item(29, 6).Value = 386
$ws.Cells.Item(30, 6).Value = 53
$ws.Cells.Item(31, 6).Value = 148
$ws.Cells.Item(32, 6).Value = 98
$ws.Cells.Item(33, 6).Value = 1217
$ws.Cells.Item(34, 6).Value = 31
$ws.Cells.Item(35, 6).Value = 47
$ws.Cells.Item(36, 6).Value = 1440
$ws.Cells.Item(37, 6).Value = 1424
$ws.Cells.Item(38, 6).Value = 916
$ws.Cells.Item(39, 6).Value = 526
$ws.Cells.Item(43, 6).Value = 1287
$ws.Cells.Item(45, 6).Value = 87
$ws.Cells.Item(46, 6).Value = 184
$ws.Cells.Item(49, 6).Value = 3757
